$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6815446615219116
$ws.Range("B1").Value = 0.6016117334365845
$ws.Range("C1").Value = 3.555123805999756
$ws.Range("D1").Value = 1.945112228393555
$ws.Range("E1").Value = 1.276267409324646
